$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("工作表3")
$ws.Activate()

# Insert a new row above the current row 2, shifting existing data down.
$ws.Rows.Item(2).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the new header row with the variable/quantity labels.
$ws.Cells.Item(2, 2).Value = "變數"
$ws.Cells.Item(2, 3).Value = "數量"

# Match the saved selection state from the edit.
$ws.Range("C3").Select()
